$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8947585821151733
$ws.Range("B1").Value = 1.762492775917053
$ws.Range("D1").Value = 1.879625201225281
$ws.Range("E1").Value = 1.112817168235779
